# Script: apply re-scrape update to france_ligue-1_2023-2024 sheet
# - Rows 25/26, 33/34/35, 42/43/44, 102/103 get re-ordered (their B:V
#   contents are relocated while the Indice column A stays sequential).
# - Two brand-new matches (Rennes-Lyon and Lens-Marseille) are appended
#   as rows 105 and 106.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1) Rows 25 <-> 26 : swap match content (columns B through V)
# ---------------------------------------------------------------
$r25 = $ws.Range("B25:V25")
$r26 = $ws.Range("B26:V26")
$v25 = $r25.Value2
$v26 = $r26.Value2
$r25.Value2 = $v26
$r26.Value2 = $v25

# ---------------------------------------------------------------
# 2) Rows 33,34,35 : rotate content (33<-35, 34<-33, 35<-34)
# ---------------------------------------------------------------
$r33 = $ws.Range("B33:V33")
$r34 = $ws.Range("B34:V34")
$r35 = $ws.Range("B35:V35")
$v33 = $r33.Value2
$v34 = $r34.Value2
$v35 = $r35.Value2
$r33.Value2 = $v35
$r34.Value2 = $v33
$r35.Value2 = $v34

# ---------------------------------------------------------------
# 3) Rows 42,43,44 : rotate content (42<-44, 43<-42, 44<-43)
# ---------------------------------------------------------------
$r42 = $ws.Range("B42:V42")
$r43 = $ws.Range("B43:V43")
$r44 = $ws.Range("B44:V44")
$v42 = $r42.Value2
$v43 = $r43.Value2
$v44 = $r44.Value2
$r42.Value2 = $v44
$r43.Value2 = $v42
$r44.Value2 = $v43

# ---------------------------------------------------------------
# 4) Rows 102 <-> 103 : swap match content (columns B through V)
# ---------------------------------------------------------------
$r102 = $ws.Range("B102:V102")
$r103 = $ws.Range("B103:V103")
$v102 = $r102.Value2
$v103 = $r103.Value2
$r102.Value2 = $v103
$r103.Value2 = $v102

# ---------------------------------------------------------------
# 5) Append two freshly-scraped matches as rows 105 and 106.
#    First clone the formatting of the last existing row (104) so the
#    new rows inherit the same styles (bold/bordered index column,
#    datetime number format for the kickoff-date column, etc.).
# ---------------------------------------------------------------
$ws.Range("A104:V104").Copy($ws.Range("A105:V105"))
$ws.Range("A104:V104").Copy($ws.Range("A106:V106"))

# Row 105 : Rennes 0 - 1 Lyon
$ws.Range("A105").Value = 104
$ws.Range("B105").Value = "france"
$ws.Range("C105").Value = "ligue-1"
$ws.Range("D105").Value = "2023-2024"
$ws.Range("E105").Value = 45242.71180555555
$ws.Range("F105").Value = "Rennes"
$ws.Range("G105").Value = 0
$ws.Range("H105").Value = "Lyon"
$ws.Range("I105").Value = 1
$ws.Range("J105").Value = 1.7
$ws.Range("K105").Value = "29/10/2023 11:02"
$ws.Range("L105").Value = 1.88
$ws.Range("M105").Value = "12/11/2023 16:18"
$ws.Range("N105").Value = 4.01
$ws.Range("O105").Value = "29/10/2023 11:02"
$ws.Range("P105").Value = 3.92
$ws.Range("Q105").Value = "12/11/2023 16:38"
$ws.Range("R105").Value = 4.38
$ws.Range("S105").Value = "29/10/2023 11:02"
$ws.Range("T105").Value = 4.09
$ws.Range("U105").Value = "12/11/2023 16:38"
$ws.Range("V105").Value = "https://www.betexplorer.com/football/france/ligue-1/rennes-lyon/lnpQ7cr0/"

# Row 106 : Lens 1 - 0 Marseille
$ws.Range("A106").Value = 105
$ws.Range("B106").Value = "france"
$ws.Range("C106").Value = "ligue-1"
$ws.Range("D106").Value = "2023-2024"
$ws.Range("E106").Value = 45242.86458333334
$ws.Range("F106").Value = "Lens"
$ws.Range("G106").Value = 1
$ws.Range("H106").Value = "Marseille"
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 2.36
$ws.Range("K106").Value = "29/10/2023 11:02"
$ws.Range("L106").Value = 2.32
$ws.Range("M106").Value = "12/11/2023 20:43"
$ws.Range("N106").Value = 3.42
$ws.Range("O106").Value = "29/10/2023 11:02"
$ws.Range("P106").Value = 3.43
$ws.Range("Q106").Value = "12/11/2023 20:42"
$ws.Range("R106").Value = 3.08
$ws.Range("S106").Value = "29/10/2023 11:02"
$ws.Range("T106").Value = 3.27
$ws.Range("U106").Value = "12/11/2023 20:43"
$ws.Range("V106").Value = "https://www.betexplorer.com/football/france/ligue-1/lens-marseille/x8QKTv5J/"
